# Daily update at 8 AM UTC
# The existing "last row" (row 56) loses its special "latest date" style
# (YYYY-MM-DD) and reverts to the normal date/time style used by every
# other row above it. A brand-new row 57 is appended with the next day's
# data and takes on the "latest date" style instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 56 (previously the last row) reverts to the standard date style.
$ws.Range("A56").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 57 with the next day's values.
$ws.Range("A57").Value = 45642
$ws.Range("B57").Value = 136
$ws.Range("C57").Value = 123
$ws.Range("D57").Value = 130

# New last row gets the "latest date" style.
$ws.Range("A57").NumberFormat = "YYYY-MM-DD"
